$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("User")

# Update the password value for the row with uid=5 (row 7: 0,1,2,3,4 nric entries precede it,
# header is row 1) from "password" to "Password".
$ws.Range("C7").Value = "Password"
